$wb = $excel.ActiveWorkbook

# Work on the 3rd sheet ("Лист3")
$ws3 = $wb.Worksheets.Item(3)

# Add the new text to cell C6 on the 3rd sheet
$ws3.Range("C6").Value = "Something on 3rd sheet"

# Select cell C6 and make this sheet the active one (tabSelected + activeTab)
$ws3.Select()
$ws3.Range("C6").Select()
